$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "url" field of the "file" sub-object (under documents[]) was removed from
# the schema table. That corresponds to row 12 of the worksheet
# (J12=url, K12=URL, L12="A URL pointing to the stored file"). Deleting the
# entire row shifts every following row up by one, which reproduces every
# other change in the diff (cascading relabelling of fields/descriptions,
# the dimension shrinking from N56 to N55, and every merged-cell range
# shifting up by one row).
$ws.Rows.Item(12).Delete()

# Row-delete correctly re-anchors the bulk of the merged cell ranges, but the
# single-cell group-header merges that used to sit at (old) rows 43 and 44
# ("Conflict of interest" / "Checklist") land, after the shift, on rows 42 and
# 43 respectively. Row 43 already keeps its merged flag from the shift, but
# row 42 needs to be (re)merged explicitly.
$ws.Range("A42").Merge()
$ws.Range("B42").Merge()
